$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# Text format first, otherwise Excel will coerce them into a number cell
# instead of keeping the original inline-string/text cell type.

$ws.Range("D2").Value = "66.592.97"
$ws.Range("E2").Value = "  -5.81%  "

$ws.Range("D3").Value = "3.207.58"
$ws.Range("E3").Value = "  -8.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.09"
$ws.Range("E5").Value = "  -5.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.80"
$ws.Range("E6").Value = "  -13.73%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.199.76"
$ws.Range("E8").Value = "  -9.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -11.51%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  -13.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.67"
$ws.Range("E11").Value = "  -9.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.496"
$ws.Range("E12").Value = "  -15.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.07"
$ws.Range("E13").Value = "  -18.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000242"
$ws.Range("E14").Value = "  -12.40%  "

$ws.Range("D15").Value = "3.719.97"
$ws.Range("E15").Value = "  -9.18%  "

$ws.Range("D16").Value = "66.592.46"
$ws.Range("E16").Value = "  -5.88%  "

$ws.Range("D17").Value = "3.209.33"
$ws.Range("E17").Value = "  -8.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "539.00"
$ws.Range("E18").Value = "  -12.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.114"
$ws.Range("E19").Value = "  -6.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.08"
$ws.Range("E20").Value = "  -16.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.98"
$ws.Range("E21").Value = "  -15.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.751"
$ws.Range("E22").Value = "  -15.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  -14.53%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.93"
$ws.Range("E24").Value = "  -13.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.29"
$ws.Range("E25").Value = "  -15.48%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.11"
$ws.Range("E27").Value = "  -17.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.98"
$ws.Range("E28").Value = "  -12.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "28.97"
$ws.Range("E29").Value = "  -14.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.54"
$ws.Range("E31").Value = "  -16.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  -13.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "531.43"
$ws.Range("E33").Value = "  -13.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.48"
$ws.Range("E34").Value = "  -20.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.66"
$ws.Range("E35").Value = "  -17.46%  "

$ws.Range("E36").Value = "  +0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.84"
$ws.Range("E37").Value = "  -7.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0419"
$ws.Range("E38").Value = "  -11.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0841"
$ws.Range("E39").Value = "  -16.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.04"
$ws.Range("E40").Value = "  -16.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.124"
$ws.Range("E41").Value = "  -14.48%  "

$ws.Range("D42").Value = "2.904.30"
$ws.Range("E42").Value = "  -13.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"

$ws.Range("D44").Value = "0.0₃0581"
$ws.Range("E44").Value = "  -21.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.259"
$ws.Range("E45").Value = "  -17.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.32"
$ws.Range("E47").Value = "  -21.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.62"
$ws.Range("E48").Value = "  -20.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.08"
$ws.Range("E49").Value = "  -19.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "123.20"
$ws.Range("E50").Value = "  -7.98%  "

$ws.Range("E51").Value = "  -14.02%  "
